$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'278.86"
$ws.Range("D3").Value = "'27.44"
$ws.Range("E3").Value = "'2.32%"
$ws.Range("D4").Value = "'4.809"
$ws.Range("E4").Value = "'2.29%"
$ws.Range("D5").Value = "'0.06333"
$ws.Range("E5").Value = "'2.35%"
$ws.Range("D6").Value = "'6.948"
$ws.Range("E6").Value = "'3.57%"
$ws.Range("D7").Value = "'3.410"
$ws.Range("E7").Value = "'7.51%"
$ws.Range("D8").Value = "'0.8797"
$ws.Range("E8").Value = "'3.44%"
$ws.Range("D9").Value = "'0.9539"
$ws.Range("E9").Value = "'4.61%"
$ws.Range("D10").Value = "'0.1473"
$ws.Range("E10").Value = "'5.05%"
$ws.Range("D11").Value = "'0.05134"
$ws.Range("E11").Value = "'0.24%"
$ws.Range("D12").Value = "'0.07346"
$ws.Range("E12").Value = "'3.55%"
$ws.Range("D13").Value = "'0.03154"
$ws.Range("E13").Value = "'1.45%"
$ws.Range("D14").Value = "'0.09066"
$ws.Range("E14").Value = "'0.32%"
$ws.Range("D15").Value = "'0.001559"
$ws.Range("E15").Value = "'0.97%"
$ws.Range("D16").Value = "'0.0006294"
$ws.Range("E16").Value = "'2.35%"
$ws.Range("D17").Value = "'0.005988"
$ws.Range("E17").Value = "'-0.01%"
$ws.Range("D18").Value = "'3.460"
$ws.Range("E18").Value = "'0.37%"
$ws.Range("E20").Value = "'2.35%"
$ws.Range("E21").Value = "'0.05%"
$ws.Range("E22").Value = "'-5.43%"
$ws.Range("D23").Value = "'0.04323"
$ws.Range("E23").Value = "'1.62%"
$ws.Range("E24").Value = "'-0.25%"
$ws.Range("E25").Value = "'6.04%"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("D27").Value = "'0.0001690"
$ws.Range("E27").Value = "'3.04%"
$ws.Range("D40").Value = "'0.04083"
$ws.Range("E40").Value = "'2.83%"
$ws.Range("D41").Value = "'0.006677"
$ws.Range("E41").Value = "'61.27%"
$ws.Range("D42").Value = "'0.1164"
$ws.Range("E42").Value = "'4.71%"
$ws.Range("D43").Value = "'0.002199"
$ws.Range("E43").Value = "'2.66%"
$ws.Range("D44").Value = "'0.01296"
$ws.Range("E44").Value = "'-2.36%"
$ws.Range("D45").Value = "'0.00005211"
$ws.Range("E45").Value = "'0.93%"
$ws.Range("E46").Value = "'-0.09%"
$ws.Range("D47").Value = "'2.377"
$ws.Range("E47").Value = "'855.16%"
$ws.Range("D48").Value = "'0.02250"
$ws.Range("E48").Value = "'-33.88%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.09%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.09%"
